# commiting POM framework for selenium
#
# The original workbook had a single sheet "LoginTestData" with one data
# row (row 2): a phone number in column A and a "Mokalpur@2021" mailto
# hyperlink in column B (row 1 is the "Username"/"Password" header row).
#
# This change:
#   1. Duplicates that data row nine more times (rows 3-11), each with the
#      same phone number and the same mailto hyperlink, formatted with the
#      existing "Hyperlink" cell style.
#   2. Moves the active selection to B2.
#   3. Adds a new, blank worksheet named "Sheet1" after "LoginTestData".
#   4. Leaves "LoginTestData" as the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginTestData")

$phone = 9594748758
$mailAddress = "mailto:Mokalpur@2021"
$mailText = "Mokalpur@2021"

for ($row = 3; $row -le 11; $row++) {
    $ws.Range("A$row").Value = $phone
    $ws.Range("B$row").Value = $mailText

    $cell = $ws.Range("B$row")
    $ws.Hyperlinks.Add($cell, $mailAddress) | Out-Null
    $cell.Style = "Hyperlink"
}

# Add the new, empty worksheet right after the existing one.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$newSheet.Name = "Sheet1"

# Keep "LoginTestData" as the active sheet, with B2 selected.
$ws.Activate()
$ws.Range("B2").Select() | Out-Null
